$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the hidden "_xlchart.v1.*" defined names left over from a
#    previously-deleted chart (workbook.xml <definedNames> block).
# ---------------------------------------------------------------------------
$nameCount = $wb.Names.Count
for ($i = $nameCount; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 2. Populate the "temp_august" column (I) for every data row (2-73) with the
#    newly-added soil temperature readings.
# ---------------------------------------------------------------------------
$augustValues = @{
    2  = 25.7;  3  = 26;    4  = 26.1;  5  = 25.8;  6  = 26.2;  7  = 26.1
    8  = 26.3;  9  = 26.3;  10 = 26.2;  11 = 26;    12 = 25.7;  13 = 26
    14 = 25.6;  15 = 25.7;  16 = 26.1;  17 = 26;    18 = 25.9;  19 = 26.2
    20 = 26.1;  21 = 26;    22 = 25.7;  23 = 26;    24 = 26.2;  25 = 26.5
    26 = 26.5;  27 = 26;    28 = 26;    29 = 26.4;  30 = 26.2;  31 = 26.6
    32 = 26.2;  33 = 26.2;  34 = 26.5;  35 = 26.2;  36 = 26.5;  37 = 26.6
    38 = 26.6;  39 = 26.5;  40 = 26.6;  41 = 26.5;  42 = 26.6;  43 = 26.6
    44 = 26.2;  45 = 26;    46 = 26;    47 = 27.1;  48 = 27.5;  49 = 27.2
    50 = 26.8;  51 = 26.9;  52 = 27;    53 = 27.3;  54 = 27.1;  55 = 27.2
    56 = 27.2;  57 = 27;    58 = 27.1;  59 = 26.9;  60 = 27.2;  61 = 27.3
    62 = 26.8;  63 = 26.6;  64 = 26.6;  65 = 26.6;  66 = 26.3;  67 = 27
    68 = 25.9;  69 = 26;    70 = 26.2;  71 = 26.1;  72 = 26.1;  73 = 25.9
}

foreach ($row in ($augustValues.Keys | Sort-Object)) {
    $ws.Range("I$row").Value = $augustValues[$row]
}

# ---------------------------------------------------------------------------
# 3. Move the active selection to I74 (just past the last data row), matching
#    the cursor position left behind after entering the new column of data.
# ---------------------------------------------------------------------------
$ws.Range("I74").Select()
